$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary numbers ---
# "VALOR MORA" total
$ws.Range("E11").Value = 50784
# "Cant. Trabajadores"
$ws.Range("C13").Value = 3
# "Cant. Periodos"
$ws.Range("F13").Value = 1

# --- Worker detail table (rows 16-19) ---
# Copy the bottom-border formatting from the last data row (19, JORGE's row)
# onto row 18 before we delete row 19, so the new last row (YORMAN, moved
# down to row 18) keeps the table's closing border.
$ws.Range("B19:J19").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 16 becomes NELSON DE JESUS CORDOBA VARGAS (previously on row 17)
$ws.Range("C16").Value = "70514964"
$ws.Range("D16").Value = "NELSON DE JESUS CORDOBA VARGAS"
$ws.Range("F16").Value = 8832

# Row 17 becomes the new worker LAURIANO CURE SUAREZ (replaces FABIO)
$ws.Range("C17").Value = "1046404907"
$ws.Range("D17").Value = "LAURIANO CURE SUAREZ"
$ws.Range("F17").Value = 8832
$ws.Range("G17").Value = 828116

# Row 18 becomes YORMAN SANTIAGO AYALA (previously on row 16)
$ws.Range("C18").Value = "20246181"
$ws.Range("D18").Value = "YORMAN SANTIAGO AYALA"
$ws.Range("E18").Value = "2001"
$ws.Range("F18").Value = 33120
$ws.Range("G18").Value = 828000

# Remove old row 19 (JORGE LUIS GRONDONA VILLEGAS) entirely; the signature
# rows below shift up automatically.
$ws.Rows("19:19").Delete()
